$d = $word.ActiveDocument

# Paragraph 1 holds the hidden **ID__...__ID** marker text.
$p = $d.Paragraphs(1)

# Add paragraph border spacing (top/left/bottom/right, 5pt "space" only,
# no visible line) and bump the left indent from 6pt (120 twips) to
# 11.25pt (225 twips).
$p.Format.Borders.DistanceFromTop = 5
$p.Format.Borders.DistanceFromLeft = 5
$p.Format.Borders.DistanceFromBottom = 5
$p.Format.Borders.DistanceFromRight = 5
$p.Format.LeftIndent = 11.25

# Replace the old topic-7 marker text with the new AF_PGI_5316_505 one,
# matching through the trailing space run so the two runs collapse into
# a single run with no stray trailing space.
$d.Content.Find.Execute("**ID__AFFARS_pgi_5316_topic_7__ID** ", $false, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_AF_PGI_5316_505__ID**", 2)
